# Add a new "Title and Content" slide at the end of the deck with some
# ideas for future development.

$p = $ppt.ActivePresentation

$newIndex = $p.Slides.Count + 1
$slide = $p.Slides.Add($newIndex, 2)   # 2 = ppLayoutText ("Title and Content")

# Title placeholder
$title = $slide.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Future development"
$title.LanguageID = "en-GB"

# Body / content placeholder - add one bullet per paragraph, setting the
# language on each new paragraph individually so every run gets tagged
# consistently (setting it once on the whole range only sticks to the
# first run).
$bullets = @(
    "Not complete",
    "Reduce necessity of internet access",
    "Improve phone compatibility",
    "Give options to add new items to shops"
)

$body = $slide.Shapes.Item(2).TextFrame.TextRange
$body.Text = $bullets[0]
$body.LanguageID = "en-GB"

for ($i = 1; $i -lt $bullets.Count; $i++) {
    $body = $slide.Shapes.Item(2).TextFrame.TextRange
    $added = $body.InsertAfter("`r" + $bullets[$i])
    $added.LanguageID = "en-GB"
}
